$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the existing "N" row (currently row 30),
# shifting it down to row 32.
$ws.Rows("30:31").Insert()

$newRange = $ws.Range("A30:L31")

# The sheet stores every value (including numeric-looking ones) as plain
# text, so force text format before assigning the values; otherwise
# strings like "0.09" would be parsed as numbers. Reset back to the
# workbook's default "Normal" style afterwards so the new rows keep the
# same (un-styled) look as the rest of the data rows.
$newRange.NumberFormat = "@"

# New row 30: Cognitive Difficulty
$ws.Range("A30").Value = "Cognitive Difficulty"
$ws.Range("B30").Value = "0.09"
$ws.Range("C30").Value = "0.15"
$ws.Range("D30").Value = "0.13"
$ws.Range("E30").Value = "0.13"
$ws.Range("F30").Value = "0.08"
$ws.Range("G30").Value = "0.07"
$ws.Range("H30").Value = "0.08"
$ws.Range("I30").Value = "0.11"
$ws.Range("J30").Value = "0.12"
$ws.Range("K30").Value = "0.08"
$ws.Range("L30").Value = "0.11"

# New row 31: Independence Difficulty
$ws.Range("A31").Value = "Independence Difficulty"
$ws.Range("B31").Value = "0.15"
$ws.Range("C31").Value = "0.21"
$ws.Range("D31").Value = "0.18"
$ws.Range("E31").Value = "0.18"
$ws.Range("F31").Value = "0.13"
$ws.Range("G31").Value = "0.1"
$ws.Range("H31").Value = "0.14"
$ws.Range("I31").Value = "0.15"
$ws.Range("J31").Value = "0.18"
$ws.Range("K31").Value = "0.12"
$ws.Range("L31").Value = "0.15"

$newRange.Style = "Normal"
